$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# New "Task" descriptions for the ICQ install rows (11-13), previously blank.
$ws.Range("D11").Value = "informasi cara install ICQ di Windows,Linux,mobile(Iphone)"
$ws.Range("D12").Value = "informasi cara install ICQ di Mac,Mobile(Symbian),Mobile(Blackberry),Mobile(Windows Phone 7)"
$ws.Range("D13").Value = "informasi cara install ICQ di Mobile(Android),Mobile(Java),Mobile(Windows Mobile),Mobile(Bada)"

# Row 12: bump daily remaining-effort values from 2 to 4 (E:M), and N12 from 0 to 1.
$ws.Range("E12:M12").Value = 4
$ws.Range("N12").Value = 1

# Row 13: bump daily remaining-effort values from 3 to 4 (E:N), and O13 from 0 to 1.
$ws.Range("E13:N13").Value = 4
$ws.Range("O13").Value = 1

# Update the active selection to match the author's last cursor position.
$ws.Range("P13").Select()
